$wb = $excel.ActiveWorkbook

# --- Sheet: Estadisticos 1P ---
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")
$ws1.Range("D2").Value = 1
$ws1.Range("E2").Value = 9
$ws1.Range("F2").Value = 29
$ws1.Range("G2").Value = 74.36
$ws1.Range("H2").Value = 6.7

$ws1.Range("D3").Value = 0
$ws1.Range("E3").Value = 1
$ws1.Range("F3").Value = 23
$ws1.Range("G3").Value = 95.83
$ws1.Range("H3").Value = 8.6

$ws1.Range("D4").Value = 0
$ws1.Range("E4").Value = 5
$ws1.Range("F4").Value = 29
$ws1.Range("G4").Value = 85.29000000000001
$ws1.Range("H4").Value = 8.1

# --- Sheet: Estadisticos 2P ---
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$ws2.Range("D2").Value = 1
$ws2.Range("E2").Value = 16
$ws2.Range("F2").Value = 22
$ws2.Range("G2").Value = 56.41
$ws2.Range("H2").Value = 6.7

$ws2.Range("D3").Value = 0
$ws2.Range("E3").Value = 2
$ws2.Range("F3").Value = 22
$ws2.Range("G3").Value = 91.67
$ws2.Range("H3").Value = 8.6

$ws2.Range("D4").Value = 0
$ws2.Range("E4").Value = 6
$ws2.Range("F4").Value = 28
$ws2.Range("G4").Value = 82.34999999999999
$ws2.Range("H4").Value = 8.1

# --- Sheet: Estadisticos Final ---
$ws3 = $wb.Worksheets.Item("Estadisticos Final")
$ws3.Range("D2").Value = 1
$ws3.Range("E2").Value = 16
$ws3.Range("F2").Value = 22
$ws3.Range("G2").Value = 56.41
$ws3.Range("H2").Value = 6.6

$ws3.Range("D3").Value = 0
$ws3.Range("E3").Value = 2
$ws3.Range("F3").Value = 22
$ws3.Range("G3").Value = 91.67
$ws3.Range("H3").Value = 8.6

$ws3.Range("D4").Value = 0
$ws3.Range("E4").Value = 6
$ws3.Range("F4").Value = 28
$ws3.Range("G4").Value = 82.34999999999999
$ws3.Range("H4").Value = 8.300000000000001

$wb.Save()
